$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.074.79'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.28%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.738.49'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +3.65%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '606.60'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +1.75%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '170.03'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +6.33%  '

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.03%  '

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +1.05%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.738.00'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +3.66%  '

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +4.82%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.370'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +5.77%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.36'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +1.72%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.88'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +3.43%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.237.13'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +3.63%  '

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +3.23%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '69.024.61'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.32%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.714.10'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +2.03%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.93'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +5.32%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '376.57'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +5.09%  '

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +5.28%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.56'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +4.03%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.01'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +6.23%  '

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +3.87%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '73.99'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.82%  '

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.16%  '

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +4.44%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.872.44'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +1.43%  '

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +3.96%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '591.18'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +6.41%  '

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.28%  '

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +5.52%  '

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +5.91%  '

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +5.94%  '

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +5.09%  '

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +3.23%  '

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.01%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '20.08'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +2.22%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '162.38'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +1.75%  '

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +4.43%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.52'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +3.82%  '

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +3.88%  '

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +1.23%  '

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.05%  '

$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = 'OKB'
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '41.19'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +2.18%  '

$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = 'BabyDogeCoin'
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0₆0311'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -1.94%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '156.38'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.02%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.98'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +5.75%  '

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +8.15%  '

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +7.29%  '
